# Add three new hint/comment strings to Sheet3, below the existing
# "What is 40.1 times 3?" / "Leeway" question (row 1) and its data row
# (row 2). Each new string lands in column B of rows 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("B3").Value = "Don't forget to include the density"
$ws.Range("B4").Value = "You may have a problem with units"
$ws.Range("B5").Value = "Did you add correctly, shithead?"

# Move/leave the active selection where the saved workbook shows it.
$ws.Range("D10").Select()
